# Enhanced AssessmentCategoryType code table to test allegations dimension
# situation for Pima: add rows 2-6 (AssessmentCategoryType 2 .. 6) to the
# AssessmentCategoryType sheet, and make that sheet the active/selected tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AssessmentCategoryType")

# Append the five new code rows below the existing "AssessmentCategoryType 1" row.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "AssessmentCategoryType 2"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "AssessmentCategoryType 3"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "AssessmentCategoryType 4"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "AssessmentCategoryType 5"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "AssessmentCategoryType 6"

# Make this sheet the active sheet/tab, with the cell below the new data selected.
$ws.Activate()
$ws.Range("B8").Select()
